{"js": "// Add a \"Chapter President\" line (italic) to the letter's signature\n// block, right before the closing \"Triangle Fraternity \u2013 University of\n// Washington\" line (and right after the blank paragraph that follows\n// \"Best,\").\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the signature paragraph: \"Triangle Fraternity \u2013 University of Washington\"\nlet signaturePara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text.trim();\n  if (text.indexOf(\"Triangle Fraternity\") === 0) {\n    signaturePara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!signaturePara) {\n  throw new Error('Could not find the \"Triangle Fraternity\" signature paragraph.');\n}\n\n// Insert a new paragraph with \"Chapter President\" just before it.\nconst newPara = signaturePara.insertParagraph(\"Chapter President\", Word.InsertLocation.before);\n\n// Make the run text italic.\nnewPara.font.italic = true;\n\n// Also italicize the paragraph mark itself (matches Word's behavior when\n// the whole paragraph, including its end-of-paragraph mark, is italicized).\nconst paraEndRange = newPara.getRange(Word.RangeLocation.end);\nparaEndRange.font.italic = true;\n\nawait context.sync();\n", "ps1": "# Add a \"Chapter President\" line (italic) to the letter's signature\n# block, right before the closing \"Triangle Fraternity - University of\n# Washington\" line (and right after the blank paragraph that follows\n# \"Best,\").\n\n$d = $word.ActiveDocument\n\n# Locate the (1-based) index of the signature paragraph that starts with\n# \"Triangle Fraternity\".\n$signatureIndex = $null\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text.TrimEnd([char]13, [char]7).StartsWith(\"Triangle Fraternity\")) {\n        $signatureIndex = $i\n        break\n    }\n}\n\nif ($null -eq $signatureIndex) {\n    throw \"Could not find the 'Triangle Fraternity' signature paragraph.\"\n}\n\n# Collapse to the very start of that paragraph and insert a new paragraph\n# (text + paragraph break) right before it.\n$signaturePara = $d.Paragraphs.Item($signatureIndex)\n$insertRange = $signaturePara.Range\n$insertRange.Collapse(1)  # wdCollapseStart\n$insertRange.InsertBefore(\"Chapter President\" + [char]13)\n\n# The newly created paragraph now occupies the same (1-based) index that\n# the signature paragraph used to have, since it was inserted right before\n# it. Re-fetch it fresh (rather than navigating via .Previous()) and\n# italicize its whole range (text + paragraph mark) so both the run and\n# the paragraph mark carry the italic formatting.\n$newPara = $d.Paragraphs.Item($signatureIndex)\n$newPara.Range.Font.Italic = 1\n"}
